$d = $word.ActiveDocument

# Locate the "KEY ACHIEVEMENTS AND IMPACT" section: the Heading2 paragraph
# with that text, through (but not including) the next Heading2 paragraph.
# We scope all Find/Replace + delete operations to this span so we don't
# touch the duplicate bullet text that also appears earlier under
# "PROFESSIONAL EXPERIENCE".

$sectionStart = 0
$sectionEnd = $d.Paragraphs.Count

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $styleName = $p.Range.ParagraphStyle.NameLocal
    if ($styleName -eq "Heading 2") {
        $txt = $p.Range.Text.Trim()
        if ($sectionStart -eq 0 -and $txt -eq "KEY ACHIEVEMENTS AND IMPACT") {
            $sectionStart = $i
        }
        elseif ($sectionStart -ne 0 -and $i -gt $sectionStart) {
            $sectionEnd = $i - 1
            break
        }
    }
}

# --- Rewrite the four "Key Achievements" bullet paragraphs to be concise,
#     impact-focused accomplishment statements. Each Find/Replace is scoped
#     to a single paragraph's Range so nothing outside the section (e.g.
#     the identical bullet text under PROFESSIONAL EXPERIENCE) is touched. ---

function Replace-InSection($oldText, $newText) {
    for ($i = $sectionStart; $i -le $sectionEnd; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -like ("*" + $oldText + "*")) {
            $p.Range.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
            return
        }
    }
}

Replace-InSection `
    "Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%" `
    "Predictive excellence: Achieved 87% voter turnout accuracy vs. 71% industry standard"

Replace-InSection `
    "Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from ±4.2% to ±2.1%" `
    "Reduced polling margins from ±4.2% to ±2.1%"

Replace-InSection `
    "Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations" `
    "Executive authority: Briefed Presidents, Congressmen, Senators, Governors on election integrity, voter sentiment and postmortem analysis"

Replace-InSection `
    "Developed longitudinal data analysis methods using geospatial techniques that improved segmentation accuracy by 34% and survey incidence rates by 28%, reducing polling costs while increasing response quality" `
    "Methodological advancement: Improved segmentation accuracy 34% and survey incidence 28%"

# --- Remove the two trailing bullet paragraphs entirely (including their
#     paragraph marks), scoped to the same section, so the "Impact" list
#     drops from six bullets to four. Walk backwards so deleting a
#     paragraph doesn't disturb indices still to be visited. ---

for ($i = $sectionEnd; $i -ge $sectionStart; $i--) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*Provided expert testimony and press briefings on electoral data integrity and demographic modeling accuracy*") {
        $p.Range.Delete()
    }
    elseif ($t -like "*Demystified FEC compliance through real-time processing systems enabling transparent campaign finance monitoring*") {
        $p.Range.Delete()
    }
}
